$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 184, pushing the existing rows 184-212 down to 185-213.
$ws.Rows.Item(184).Insert()

# Populate the newly inserted row 184 with the new weekly price entry.
$ws.Cells.Item(184, 1).Value = 8
$ws.Cells.Item(184, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(184, 3).Value = "Coquimbo"
$ws.Cells.Item(184, 4).Value = 44748
$ws.Cells.Item(184, 5).Value = 4
$ws.Cells.Item(184, 6).Value = 100112037
$ws.Cells.Item(184, 7).Value = "Cebollín"
$ws.Cells.Item(184, 8).Value = "Sin especificar"
$ws.Cells.Item(184, 9).Value = "Primera"
$ws.Cells.Item(184, 10).Value = 1400
$ws.Cells.Item(184, 11).Value = 1400
$ws.Cells.Item(184, 12).Value = 1600
$ws.Cells.Item(184, 13).Value = 1500
$ws.Cells.Item(184, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(184, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(184, 16).Value = 250
$ws.Cells.Item(184, 17).Value = 6
$ws.Cells.Item(184, 18).Value = "Hortaliza"
